$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds daily data rows 2..328, with dates (Excel serials)
# in column A (formatted as date, same style as the cell above) and 0 values
# in columns B, C, D. We extend the data through row 343 (2021-08-09),
# continuing the existing pattern: date in A incremented by 1 day each row,
# 0 in B, C, D.

$lastRow = 328
$startRow = $lastRow + 1
$endRow = 343
$startSerial = 44403

# Carry the date-column formatting (number format, alignment, font, border)
# down from the last existing row onto the new rows before filling values.
$ws.Range("A" + $lastRow).Copy() | Out-Null
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt ($endRow - $startRow + 1); $i++) {
    $r = $startRow + $i
    $serial = $startSerial + $i

    $ws.Cells.Item($r, 1).Value = $serial
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
